# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the 0ac9bb0a-2e1d-43d6-b72f-ee44e997326d.md handoff batch, and marks those
# rows' Priority column as "ht" (handoff type) on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Rows 7, 8, 9, 10, 11, 13 correspond to the files that were part of this
# handoff batch (0ac9bb0a, 101e14b9, 176599fd, 24577130, 5ada6771, 826a9881).
$rows = @(7, 8, 9, 10, 11, 13)

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" column G
    $overview.Range("G$r").Value = "2016-08-17 14:19:34"

    # zh-cn sheet: "Latest Handoff Datetime" column H, and Priority column E
    $zhcn.Range("H$r").Value = "2016-08-17 14:19:29"
    $zhcn.Range("E$r").Value = "ht"

    # de-de sheet: "Latest Handoff Datetime" column H, and Priority column E
    $dede.Range("H$r").Value = "2016-08-17 14:19:34"
    $dede.Range("E$r").Value = "ht"
}
